# Revert "Fixed if statement":
# Removes the "Succinate_export" row (with its flux value) that the
# previous commit had introduced, shifting all subsequent rows back up
# by one and shrinking the shared-strings table accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(173).Delete()
